$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 3
$ws.Range("E3").Value = 4
$ws.Range("G3").Value = -3
$ws.Range("H3").Value = 13

# Update the selected cell to B3
$ws.Range("B3").Select()
